$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g3.3")

# Update period labels in column A
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "2024/2010"
}
for ($r = 8; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = "2024/2023"
}

# Update numeric values in columns C, D, E
$values = @{
    2  = @{ C = 13.66779441832582;  D = 26.7235526005577;   E = 20.66221432948787 }
    3  = @{ C = 3.415857052623439;  D = 4.136794125280074;  E = -0.2965810792255885 }
    4  = @{ C = 31.65866977514726;  D = 47.31588066856978;  E = 24.32964737999714 }
    5  = @{ C = 27.62931238300965;  D = 56.4060786169926;   E = -8.314600061897027 }
    6  = @{ C = 25.7859331629938;   D = 63.05668613763338;  E = 27.00653331278087 }
    7  = @{ C = 12.72700105075479; D = 6.947868756221642;  E = 71.29214324958792 }
    8  = @{ C = -0.1844570117514044; D = 3.344993832340304;  E = 4.207648837716005 }
    9  = @{ C = -1.646545063193328; D = 2.941735470000539;  E = 1.012714575856211 }
    10 = @{ C = 6.823012050586064;  D = 9.607515568198988;  E = 20.68178279931152 }
    11 = @{ C = 1.726868339108867;  D = 5.621579248202657;  E = 4.988522370781978 }
    12 = @{ C = 0.2860830937304382; D = 3.530685640615538;  E = 2.882952973720609 }
    13 = @{ C = 1.833274588282219;  D = 3.465999423009913;  E = 13.67806947759216 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row].C
    $ws.Cells.Item($row, 4).Value = $values[$row].D
    $ws.Cells.Item($row, 5).Value = $values[$row].E
}
